# Update the cryptocurrency price / volume table with freshly scraped values.
# Rows 17/18 and 43/44 also have their Coin name + Link swapped (re-ranked).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.392.06'
$ws.Range('E2').Value = '  +1.51%  '
$ws.Range('D3').Value = '3.374.58'
$ws.Range('E3').Value = '  +0.81%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '590.70'
$ws.Range('E5').Value = '  +5.69%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '186.49'
$ws.Range('E6').Value = '  -1.58%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.591'
$ws.Range('E8').Value = '  +1.16%  '
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.586'
$ws.Range('E10').Value = '  +0.30%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '47.33'
$ws.Range('E11').Value = '  +1.13%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000273'
$ws.Range('E12').Value = '  +1.04%  '
$ws.Range('D13').Value = '3.921.60'
$ws.Range('E13').Value = '  +1.00%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '635.70'
$ws.Range('E14').Value = '  +6.04%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.58'
$ws.Range('E15').Value = '  -0.48%  '
$ws.Range('D16').Value = '67.579.41'
$ws.Range('E16').Value = '  +1.67%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.384.40'
$ws.Range('E17').Value = '  +1.11%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.119'
$ws.Range('E18').Value = '  +0.90%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.99'
$ws.Range('E19').Value = '  -0.13%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.14'
$ws.Range('E20').Value = '  +0.77%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.909'
$ws.Range('E21').Value = '  +0.70%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '17.98'
$ws.Range('E22').Value = '  -3.45%  '
$ws.Range('E23').Value = '  +1.35%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '99.12'
$ws.Range('E24').Value = '  -0.64%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.01'
$ws.Range('E25').Value = '  +0.80%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.85'
$ws.Range('E26').Value = '  +3.52%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.64'
$ws.Range('E27').Value = '  +0.67%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '32.47'
$ws.Range('E28').Value = '  +6.08%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.65'
$ws.Range('E29').Value = '  -0.58%  '
$ws.Range('E30').Value = '  +2.90%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '606.69'
$ws.Range('E31').Value = '  +4.36%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.77'
$ws.Range('E32').Value = '  -4.22%  '
$ws.Range('D33').Value = '4.007.43'
$ws.Range('E33').Value = '  +7.99%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.07'
$ws.Range('E34').Value = '  +0.44%  '
$ws.Range('E35').Value = '  +0.84%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('E36').Value = '  +0.03%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '56.11'
$ws.Range('E37').Value = '  -0.63%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.85'
$ws.Range('E38').Value = '  +6.75%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.131'
$ws.Range('E39').Value = '  +3.30%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '33.62'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.23'
$ws.Range('E41').Value = '  +0.93%  '
$ws.Range('D42').Value = '0.0₃0700'
$ws.Range('E42').Value = '  -0.96%  '
$ws.Range('B43').Value = 'ApeXProtocol'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.39'
$ws.Range('E43').Value = '  +0.32%  '
$ws.Range('B44').Value = 'TheGraph'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.341'
$ws.Range('E44').Value = '  -0.20%  '
$ws.Range('E45').Value = '  +0.48%  '
$ws.Range('E46').Value = '  +0.17%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.58'
$ws.Range('E47').Value = '  +0.25%  '
$ws.Range('E48').Value = '  +0.33%  '
$ws.Range('E49').Value = '  +9.60%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.86'
$ws.Range('E50').Value = '  -18.79%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '128.12'
$ws.Range('E51').Value = '  +2.97%  '
